$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 670.7213
$ws.Range("J17").Value = 670.7213
$ws.Range("L17").Value = 2012.1639
$ws.Range("N17").Value = -2348.1639
$ws.Range("H28").Value = 476.47058
$ws.Range("I28").Value = 521.3333
$ws.Range("J28").Value = 140
$ws.Range("K28").Value = 521.3333
$ws.Range("L28").Value = 140
$ws.Range("M28").Value = -36.33330000000001
$ws.Range("N28").Value = -1110
$ws.Range("H95").Value = 20599.285
$ws.Range("J95").Value = 20599.285
$ws.Range("L95").Value = 20599.285
$ws.Range("N95").Value = -26091.285
$ws.Range("H98").Value = 2513.6667
$ws.Range("I98").Value = 2049.303
$ws.Range("K98").Value = 2049.303
$ws.Range("M98").Value = -551.3029999999999
$ws.Range("H122").Value = 2513.6667
$ws.Range("I122").Value = 2049.303
$ws.Range("K122").Value = 6147.909
$ws.Range("M122").Value = -3697.909
$ws.Range("H132").Value = 3902.3044
$ws.Range("I132").Value = 3430.8333
$ws.Range("K132").Value = 10292.4999
$ws.Range("M132").Value = -7762.499899999999
$ws.Range("H137").Value = 2634.4
$ws.Range("I137").Value = 2027.1818
$ws.Range("J137").Value = 3111.5
$ws.Range("K137").Value = 6081.5454
$ws.Range("L137").Value = 9334.5
$ws.Range("M137").Value = -3531.5454
$ws.Range("N137").Value = -14434.5
$ws.Range("H138").Value = 9516.277
$ws.Range("I138").Value = 7469.1816
$ws.Range("J138").Value = 10039.953
$ws.Range("K138").Value = 22407.5448
$ws.Range("L138").Value = 30119.859
$ws.Range("M138").Value = -17267.5448
$ws.Range("N138").Value = -40399.859

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H32").Value = 1306.6274
$ws.Range("I32").Value = 823.65216
$ws.Range("K32").Value = 823.65216
$ws.Range("M32").Value = -536.65216
$ws.Range("H61").Value = 12222.1
$ws.Range("I61").Value = 8162
$ws.Range("J61").Value = 18312.25
$ws.Range("K61").Value = 8162
$ws.Range("L61").Value = 18312.25
$ws.Range("M61").Value = -7950
$ws.Range("N61").Value = -18736.25
$ws.Range("H62").Value = 10000
$ws.Range("J62").Value = 10000
$ws.Range("L62").Value = 10000
$ws.Range("N62").Value = -11248
$ws.Range("H65").Value = 10000
$ws.Range("J65").Value = 10000
$ws.Range("L65").Value = 30000
$ws.Range("N65").Value = -36240
$ws.Range("H74").Value = 3163.55
$ws.Range("I74").Value = 2294.1
$ws.Range("J74").Value = 4033
$ws.Range("K74").Value = 2294.1
$ws.Range("L74").Value = 4033
$ws.Range("M74").Value = -1420.1
$ws.Range("N74").Value = -5781
$ws.Range("H77").Value = 3163.55
$ws.Range("I77").Value = 2294.1
$ws.Range("J77").Value = 4033
$ws.Range("K77").Value = 11470.5
$ws.Range("L77").Value = 20165
$ws.Range("M77").Value = -7102.5
$ws.Range("N77").Value = -28901
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H104").Value = 39998.75
$ws.Range("J104").Value = 39998.75
$ws.Range("L104").Value = 39998.75
$ws.Range("N104").Value = -46986.75
$ws.Range("H136").Value = 12222.1
$ws.Range("I136").Value = 8162
$ws.Range("J136").Value = 18312.25
$ws.Range("K136").Value = 24486
$ws.Range("L136").Value = 54936.75
$ws.Range("M136").Value = -21936
$ws.Range("N136").Value = -60036.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 15605.723
$ws.Range("I134").Value = 16478.46
$ws.Range("J134").Value = 13336.6
$ws.Range("K134").Value = 49435.38
$ws.Range("L134").Value = 40009.8
$ws.Range("M134").Value = -46900.38
$ws.Range("N134").Value = -45079.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3053.2205
$ws.Range("J31").Value = 3043.3958
$ws.Range("L31").Value = 3043.3958
$ws.Range("N31").Value = -3633.3958
$ws.Range("H34").Value = 3053.2205
$ws.Range("J34").Value = 3043.3958
$ws.Range("L34").Value = 3043.3958
$ws.Range("N34").Value = -3447.3958
$ws.Range("H58").Value = 6845.769
$ws.Range("I58").Value = 4316.1665
$ws.Range("J58").Value = 9014
$ws.Range("K58").Value = 4316.1665
$ws.Range("L58").Value = 9014
$ws.Range("M58").Value = -4113.1665
$ws.Range("N58").Value = -9420
$ws.Range("H86").Value = 3366
$ws.Range("I86").Value = 3139.2
$ws.Range("J86").Value = 4500
$ws.Range("K86").Value = 3139.2
$ws.Range("L86").Value = 4500
$ws.Range("M86").Value = -2016.2
$ws.Range("N86").Value = -6746
$ws.Range("H89").Value = 3366
$ws.Range("I89").Value = 3139.2
$ws.Range("J89").Value = 4500
$ws.Range("K89").Value = 15696
$ws.Range("L89").Value = 22500
$ws.Range("M89").Value = -10080
$ws.Range("N89").Value = -33732
$ws.Range("H106").Value = 56990
$ws.Range("J106").Value = 56990
$ws.Range("L106").Value = 56990
$ws.Range("N106").Value = -59514
$ws.Range("H122").Value = 2799.4
$ws.Range("I122").Value = 2719.3
$ws.Range("J122").Value = 2959.6
$ws.Range("K122").Value = 8157.900000000001
$ws.Range("L122").Value = 8878.799999999999
$ws.Range("M122").Value = -5707.900000000001
$ws.Range("N122").Value = -13778.8
$ws.Range("H132").Value = 2792.2666
$ws.Range("I132").Value = 2792.2666
$ws.Range("K132").Value = 8376.799800000001
$ws.Range("M132").Value = -5846.799800000001
$ws.Range("H134").Value = 4886.2964
$ws.Range("I134").Value = 4357.2
$ws.Range("J134").Value = 11500
$ws.Range("K134").Value = 13071.6
$ws.Range("L134").Value = 34500
$ws.Range("M134").Value = -10536.6
$ws.Range("N134").Value = -39570
$ws.Range("H136").Value = 6845.769
$ws.Range("I136").Value = 4316.1665
$ws.Range("J136").Value = 9014
$ws.Range("K136").Value = 12948.4995
$ws.Range("L136").Value = 27042
$ws.Range("M136").Value = -10398.4995
$ws.Range("N136").Value = -32142

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 36799720
$ws.Range("I4").Value = 41155676
$ws.Range("J4").Value = 500098
$ws.Range("K4").Value = 123467028
$ws.Range("L4").Value = 1500294
$ws.Range("M4").Value = -123466916
$ws.Range("N4").Value = -1500518
$ws.Range("H34").Value = 3370
$ws.Range("J34").Value = 5525
$ws.Range("L34").Value = 16575
$ws.Range("N34").Value = -16743
$ws.Range("H122").Value = 1567
$ws.Range("J122").Value = 1664.5
$ws.Range("L122").Value = 14980.5
$ws.Range("N122").Value = -19880.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H105").Value = 27647.2
$ws.Range("J105").Value = 27647.2
$ws.Range("L105").Value = 27647.2
$ws.Range("N105").Value = -34635.2
$ws.Range("H132").Value = 5917.5557
$ws.Range("I132").Value = 5917.5557
$ws.Range("K132").Value = 17752.6671
$ws.Range("M132").Value = -15222.6671
$ws.Range("H135").Value = 127000
$ws.Range("J135").Value = 127000
$ws.Range("L135").Value = 127000
$ws.Range("N135").Value = -137140

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 10496.75
$ws.Range("J101").Value = 10496.75
$ws.Range("L101").Value = 10496.75
$ws.Range("N101").Value = -16986.75
$ws.Range("H136").Value = 2081.4546
$ws.Range("I136").Value = 1288.2858
$ws.Range("K136").Value = 3864.8574
$ws.Range("M136").Value = -1314.8574
$ws.Range("H141").Value = 83331.664
$ws.Range("J141").Value = 83331.664
$ws.Range("L141").Value = 83331.664
$ws.Range("N141").Value = -93691.664

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H132").Value = 8867
$ws.Range("I132").Value = 11121.286
$ws.Range("J132").Value = 5711
$ws.Range("K132").Value = 33363.858
$ws.Range("L132").Value = 17133
$ws.Range("M132").Value = -30833.858
$ws.Range("N132").Value = -22193
$ws.Range("H136").Value = 6721.784
$ws.Range("I136").Value = 6596.628
$ws.Range("J136").Value = 7394.5
$ws.Range("K136").Value = 19789.884
$ws.Range("L136").Value = 22183.5
$ws.Range("M136").Value = -17239.884
$ws.Range("N136").Value = -27283.5
